$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 105277
$ws.Range("B3").Value = 56035
$ws.Range("G3").Value = 1202
$ws.Range("B4").Value = 41876
$ws.Range("B5").Value = 21759
$ws.Range("B6").Value = 51842
$ws.Range("H6").Value = 2005
$ws.Range("I6").Value = 0.828428927680798
$ws.Range("B7").Value = 28566
$ws.Range("B8").Value = 34824
$ws.Range("B9").Value = 58191
$ws.Range("B10").Value = 54043
$ws.Range("B11").Value = 36142
$ws.Range("B13").Value = 27612
$ws.Range("C14").Value = 229783
$ws.Range("B15").Value = 18726
$ws.Range("B16").Value = 20096
$ws.Range("B17").Value = 17120
$ws.Range("B18").Value = 12110
$ws.Range("D18").Value = 85
$ws.Range("G18").Value = 46
$ws.Range("H18").Value = 77
$ws.Range("I18").Value = 0.8311688311688312
$ws.Range("B20").Value = 9910
$ws.Range("B22").Value = 14044
$ws.Range("B23").Value = 18615
$ws.Range("B24").Value = 12303
$ws.Range("B25").Value = 23750
$ws.Range("B26").Value = 12172
$ws.Range("B27").Value = 36002
$ws.Range("B28").Value = 24700
$ws.Range("B29").Value = 23092
$ws.Range("B30").Value = 13386
$ws.Range("B33").Value = 18580
$ws.Range("H33").Value = 1165
$ws.Range("I33").Value = 0.8652360515021459
$ws.Range("B34").Value = 7142
$ws.Range("B36").Value = 13267
$ws.Range("C36").Value = 144
$ws.Range("B37").Value = 29950
$ws.Range("B38").Value = 19233
$ws.Range("B39").Value = 17626
$ws.Range("H39").Value = 440
$ws.Range("I39").Value = 0.8863636363636364
$ws.Range("B40").Value = 13919
$ws.Range("C41").Value = 128958
$ws.Range("B42").Value = 19029
$ws.Range("B43").Value = 15213
$ws.Range("B44").Value = 12496
$ws.Range("B46").Value = 9613
$ws.Range("B47").Value = 9461
$ws.Range("B49").Value = 12537
$ws.Range("H49").Value = 592
$ws.Range("I49").Value = 0.9003378378378378
$ws.Range("B50").Value = 9731
$ws.Range("B51").Value = 16169
$ws.Range("B52").Value = 14537
$ws.Range("F52").Value = 488
$ws.Range("H52").Value = 506
$ws.Range("I52").Value = 0.9644268774703557
$ws.Range("B53").Value = 13203
$ws.Range("B55").Value = 14405
$ws.Range("B56").Value = 15384
$ws.Range("B57").Value = 11463
$ws.Range("B58").Value = 24431
$ws.Range("B59").Value = 15344
$ws.Range("B60").Value = 4719
$ws.Range("B61").Value = 7506
$ws.Range("B62").Value = 13420
$ws.Range("B63").Value = 15590
$ws.Range("B64").Value = 13195
$ws.Range("B66").Value = 18218
$ws.Range("B68").Value = 8149
$ws.Range("B69").Value = 11250
$ws.Range("D69").Value = 548
$ws.Range("B71").Value = 12324
$ws.Range("B72").Value = 9774
$ws.Range("B73").Value = 9473
$ws.Range("B76").Value = 4312
$ws.Range("B77").Value = 12674
$ws.Range("D77").Value = 98
$ws.Range("F77").Value = 265
$ws.Range("H77").Value = 292
$ws.Range("I77").Value = 0.9075342465753424
$ws.Range("B78").Value = 9683
$ws.Range("B79").Value = 7286
$ws.Range("B80").Value = 9775
$ws.Range("B81").Value = 5060
$ws.Range("B83").Value = 8387
$ws.Range("B84").Value = 7357
$ws.Range("B87").Value = 5261
$ws.Range("H87").Value = 119
$ws.Range("I87").Value = 0.6974789915966386
$ws.Range("B89").Value = 13053
$ws.Range("D90").Value = 282
$ws.Range("B92").Value = 7710
$ws.Range("B93").Value = 9699
$ws.Range("B94").Value = 9856
$ws.Range("B95").Value = 12092
$ws.Range("B96").Value = 5643
$ws.Range("B97").Value = 7337
$ws.Range("B98").Value = 8011